# Updates cryptos list D (Price) / E (Volume 1h) columns for rows 2-51
# to match the refreshed figures from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.581.30"
$ws.Range("E2").Value = "  -2.43%  "

$ws.Range("D3").Value = "2.371.94"
$ws.Range("E3").Value = "  -3.98%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'311.03"
$ws.Range("E5").Value = "  -2.39%  "

$ws.Range("D6").Value = "'86.32"
$ws.Range("E6").Value = "  -6.52%  "

$ws.Range("E7").Value = "  -4.18%  "

$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("D10").Value = "'0.0838"
$ws.Range("E10").Value = "  -2.96%  "

$ws.Range("E11").Value = "  -8.50%  "

$ws.Range("D12").Value = "'0.110"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").Value = "2.733.63"
$ws.Range("E13").Value = "  -4.12%  "

$ws.Range("D14").Value = "'6.54"
$ws.Range("E14").Value = "  -4.91%  "

$ws.Range("D15").Value = "'15.09"
$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").Value = "2.374.93"
$ws.Range("E16").Value = "  -4.40%  "

$ws.Range("D17").Value = "'0.757"
$ws.Range("E17").Value = "  -4.35%  "

$ws.Range("D18").Value = "40.515.24"
$ws.Range("E18").Value = "  -2.52%  "

$ws.Range("D19").Value = "0.0₃0911"
$ws.Range("E19").Value = "  -3.41%  "

$ws.Range("D20").Value = "'6.13"
$ws.Range("E20").Value = "  -4.81%  "

$ws.Range("D21").Value = "'68.62"
$ws.Range("E21").Value = "  -2.97%  "

$ws.Range("D22").Value = "'10.76"
$ws.Range("E22").Value = "  -4.69%  "

$ws.Range("D23").Value = "'235.28"
$ws.Range("E23").Value = "  -2.13%  "

$ws.Range("D24").Value = "'2.58"
$ws.Range("E24").Value = "  -6.18%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("E26").Value = "  -8.29%  "

$ws.Range("D27").Value = "'23.89"
$ws.Range("E27").Value = "  -3.40%  "

$ws.Range("E28").Value = "  -2.22%  "

$ws.Range("D29").Value = "'9.27"
$ws.Range("E29").Value = "  -4.33%  "

$ws.Range("D30").Value = "'34.17"
$ws.Range("E30").Value = "  -6.22%  "

$ws.Range("D31").Value = "'154.38"
$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").Value = "'5.20"
$ws.Range("E33").Value = "  -4.49%  "

$ws.Range("D34").Value = "'0.0728"
$ws.Range("E34").Value = "  -4.54%  "

$ws.Range("E35").Value = "  -5.77%  "

$ws.Range("E36").Value = "  -2.44%  "

$ws.Range("D37").Value = "'16.20"
$ws.Range("E37").Value = "  -6.02%  "

$ws.Range("E38").Value = "  -3.83%  "

$ws.Range("D39").Value = "'2.75"
$ws.Range("E39").Value = "  -5.15%  "

$ws.Range("E40").Value = "  -7.83%  "

$ws.Range("D41").Value = "'3.84"
$ws.Range("E41").Value = "  -3.58%  "

$ws.Range("D42").Value = "'2.39"
$ws.Range("E42").Value = "  -3.84%  "

$ws.Range("D43").Value = "1.958.33"
$ws.Range("E43").Value = "  -1.47%  "

$ws.Range("E44").Value = "  -5.18%  "

$ws.Range("D45").Value = "'17.85"
$ws.Range("E45").Value = "  -5.24%  "

$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("E47").Value = "  -9.16%  "

$ws.Range("D48").Value = "2.596.09"
$ws.Range("E48").Value = "  -4.13%  "

$ws.Range("D49").Value = "'92.99"
$ws.Range("E49").Value = "  -4.58%  "

$ws.Range("D50").Value = "'71.89"
$ws.Range("E50").Value = "  -5.29%  "

$ws.Range("D51").Value = "'50.10"
$ws.Range("E51").Value = "  -4.24%  "
